$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 30, shifting the existing rows 30-32 down to 31-33.
$ws.Rows.Item(30).Insert()

# Populate the newly inserted row 30 with a copy of the surrounding row's
# constant values plus the new record's data.
$ws.Range("A30").Value = 8
$ws.Range("B30").Value = "Terminal La Palmera de La Serena"
$ws.Range("C30").Value = "Coquimbo"
$ws.Range("D30").Value = 44449
$ws.Range("D30").NumberFormat = $ws.Range("D31").NumberFormat
$ws.Range("E30").Value = 4
$ws.Range("F30").Value = 100112052
$ws.Range("G30").Value = "Albahaca"
$ws.Range("H30").Value = "Sin especificar"
$ws.Range("I30").Value = "Primera"
$ws.Range("J30").Value = 700
$ws.Range("K30").Value = 4000
$ws.Range("L30").Value = 4500
$ws.Range("M30").Value = 4250
$ws.Range("N30").Value = "$/paquete"
$ws.Range("O30").Value = "Región de Arica y Parinacota"
$ws.Range("P30").Value = 4250
$ws.Range("Q30").Value = 1
$ws.Range("R30").Value = "Hortaliza"
